$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.84264716618951
$ws.Range("C2").Value = 7.012486002881392
$ws.Range("D2").Value = 14.58924615032959
$ws.Range("E2").Value = 15.53506472457622
$ws.Range("G2").Value = 56.93001939872304
$ws.Range("H2").Value = 21.24060571415182
$ws.Range("J2").Value = 9.108652014879013
$ws.Range("K2").Value = 14.39671618033581
$ws.Range("M2").Value = 19.20611913869215
$ws.Range("N2").Value = 23.52835691381003

$ws.Range("B3").Value = 14.67589437981519
$ws.Range("C3").Value = 6.904718026049555
$ws.Range("D3").Value = 14.56765974190672
$ws.Range("E3").Value = 15.53580892268196
$ws.Range("G3").Value = 56.81390448312573
$ws.Range("H3").Value = 21.26364709158396
$ws.Range("J3").Value = 9.125432257104261
$ws.Range("K3").Value = 14.2966180340384
$ws.Range("M3").Value = 19.17703502480393
$ws.Range("N3").Value = 23.5776912050726

$ws.Range("B4").Value = 14.57663346378183
$ws.Range("C4").Value = 6.839799411768619
$ws.Range("D4").Value = 14.55732367638774
$ws.Range("E4").Value = 15.53901056522919
$ws.Range("G4").Value = 56.75494560857468
$ws.Range("H4").Value = 21.28105068793413
$ws.Range("J4").Value = 9.13656237759244
$ws.Range("K4").Value = 14.23840569250702
$ws.Range("M4").Value = 19.162982256956
$ws.Range("N4").Value = 23.60989521542459

$ws.Range("B5").Value = 14.53701717093531
$ws.Range("C5").Value = 6.81369535211315
$ws.Range("D5").Value = 14.5538482822083
$ws.Range("E5").Value = 15.54100606418148
$ws.Range("G5").Value = 56.73403197821754
$ws.Range("H5").Value = 21.28896089511353
$ws.Range("J5").Value = 9.141306291379282
$ws.Range("K5").Value = 14.21552253267614
$ws.Range("M5").Value = 19.15821676455022
$ws.Range("N5").Value = 23.62349991679551

$ws.Range("B6").Value = 14.53049059432606
$ws.Range("C6").Value = 6.809383025242255
$ws.Range("D6").Value = 14.55331575964625
$ws.Range("E6").Value = 15.54137915091221
$ws.Range("G6").Value = 56.73074753194829
$ws.Range("H6").Value = 21.29032376754155
$ws.Range("J6").Value = 9.142106605628744
$ws.Range("K6").Value = 14.21177407989568
$ws.Range("M6").Value = 19.15748362032739
$ws.Range("N6").Value = 23.62578804812655

$ws.Range("B7").Value = 14.5760957493074
$ws.Range("C7").Value = 6.83944589558504
$ws.Range("D7").Value = 14.55727382001542
$ws.Range("E7").Value = 15.53903467957022
$ws.Range("G7").Value = 56.7546509455702
$ws.Range("H7").Value = 21.28115405635542
$ws.Range("J7").Value = 9.136625511818194
$ws.Range("K7").Value = 14.2380936577978
$ws.Range("M7").Value = 19.1629140911281
$ws.Range("N7").Value = 23.61007674392972

$ws.Range("B8").Value = 14.78453089472474
$ws.Range("C8").Value = 6.975089717734106
$ws.Range("D8").Value = 14.58119938532268
$ws.Range("E8").Value = 15.5347521149786
$ws.Range("G8").Value = 56.88742976767858
$ws.Range("H8").Value = 21.24787421529147
$ws.Range("J8").Value = 9.114266410720164
$ws.Range("K8").Value = 14.36154180580531
$ws.Range("M8").Value = 19.19530401087388
$ws.Range("N8").Value = 23.54497065287661

$ws.Range("B9").Value = 15.21592092888974
$ws.Range("C9").Value = 7.249409463736218
$ws.Range("D9").Value = 14.65112306519546
$ws.Range("E9").Value = 15.54809486381247
$ws.Range("G9").Value = 57.24510854807217
$ws.Range("H9").Value = 21.20847429182828
$ws.Range("J9").Value = 9.076966649255834
$ws.Range("K9").Value = 14.62836209544508
$ws.Range("M9").Value = 19.28879253881139
$ws.Range("N9").Value = 23.43245581229192

$ws.Range("B10").Value = 15.54370776460465
$ws.Range("C10").Value = 7.453874385711136
$ws.Range("D10").Value = 14.71629059944079
$ws.Range("E10").Value = 15.57108747793728
$ws.Range("G10").Value = 57.56626240343443
$ws.Range("H10").Value = 21.19532288988929
$ws.Range("J10").Value = 9.053532511936215
$ws.Range("K10").Value = 14.83800342233925
$ws.Range("M10").Value = 19.37541522021789
$ws.Range("N10").Value = 23.35900771952594

$ws.Range("B11").Value = 15.69453992330536
$ws.Range("C11").Value = 7.547077359161925
$ws.Range("D11").Value = 14.74887210630864
$ws.Range("E11").Value = 15.58439179060009
$ws.Range("G11").Value = 57.72477762026538
$ws.Range("H11").Value = 21.19277293551779
$ws.Range("J11").Value = 9.043729367217253
$ws.Range("K11").Value = 14.9359944373931
$ws.Range("M11").Value = 19.41863108201653
$ws.Range("N11").Value = 23.32759074933146

$ws.Range("B12").Value = 15.7518498579188
$ws.Range("C12").Value = 7.582362120299227
$ws.Range("D12").Value = 14.76162613358571
$ws.Range("E12").Value = 15.58983664625088
$ws.Range("G12").Value = 57.78656275178587
$ws.Range("H12").Value = 21.19230082618831
$ws.Range("J12").Value = 9.040140085958548
$ws.Range("K12").Value = 14.97344909207889
$ws.Range("M12").Value = 19.43553537636688
$ws.Range("N12").Value = 23.31598054161286

$ws.Range("B13").Value = 15.73949940997371
$ws.Range("C13").Value = 7.574763894485876
$ws.Range("D13").Value = 14.7588609219976
$ws.Range("E13").Value = 15.58864594625755
$ws.Range("G13").Value = 57.7731784538753
$ws.Range("H13").Value = 21.19238055851205
$ws.Range("J13").Value = 9.040907638153339
$ws.Range("K13").Value = 14.96536758277651
$ws.Range("M13").Value = 19.43187088877915
$ws.Range("N13").Value = 23.31846825796898

$ws.Range("B14").Value = 15.69925128581635
$ws.Range("C14").Value = 7.549980647278886
$ws.Range("D14").Value = 14.7499130852209
$ws.Range("E14").Value = 15.58483160344021
$ws.Range("G14").Value = 57.72982563535061
$ws.Range("H14").Value = 21.19272420559034
$ws.Range("J14").Value = 9.043431612374397
$ws.Range("K14").Value = 14.93906906069004
$ws.Range("M14").Value = 19.42001105001323
$ws.Range("N14").Value = 23.3266298251246

$ws.Range("B15").Value = 15.67462168841471
$ws.Range("C15").Value = 7.534797930135884
$ws.Range("D15").Value = 14.74448627645262
$ws.Range("E15").Value = 15.58254811273748
$ws.Range("G15").Value = 57.70349898465349
$ws.Range("H15").Value = 21.19299896140356
$ws.Range("J15").Value = 9.044993623137227
$ws.Range("K15").Value = 14.92300482950249
$ws.Range("M15").Value = 19.41281653367001
$ws.Range("N15").Value = 23.33166635549427

$ws.Range("B16").Value = 15.53388049154863
$ws.Range("C16").Value = 7.44778403536746
$ws.Range("D16").Value = 14.71421989743662
$ws.Range("E16").Value = 15.57027506728362
$ws.Range("G16").Value = 57.5561507310199
$ws.Range("H16").Value = 21.19555860449057
$ws.Range("J16").Value = 9.054190392950131
$ws.Range("K16").Value = 14.83164975223761
$ws.Range("M16").Value = 19.37266692469031
$ws.Range("N16").Value = 23.36110102216288

$ws.Range("B17").Value = 15.44794215198037
$ws.Range("C17").Value = 7.394426293062881
$ws.Range("D17").Value = 14.69640035044955
$ws.Range("E17").Value = 15.56347305953756
$ws.Range("G17").Value = 57.46892006828088
$ws.Range("H17").Value = 21.19800803557905
$ws.Range("J17").Value = 9.060051637874764
$ws.Range("K17").Value = 14.77625748749874
$ws.Range("M17").Value = 19.34900656510361
$ws.Range("N17").Value = 23.37966907330682

$ws.Range("B18").Value = 15.39867784376897
$ws.Range("C18").Value = 7.363756420555171
$ws.Range("D18").Value = 14.68642781407488
$ws.Range("E18").Value = 15.55982865247487
$ws.Range("G18").Value = 57.41991819542688
$ws.Range("H18").Value = 21.19974000260863
$ws.Range("J18").Value = 9.063503569032651
$ws.Range("K18").Value = 14.74464614071013
$ws.Range("M18").Value = 19.33575699751832
$ws.Range("N18").Value = 23.39053670574276

$ws.Range("B19").Value = 15.38202782276191
$ws.Range("C19").Value = 7.353376713844101
$ws.Range("D19").Value = 14.68309900312296
$ws.Range("E19").Value = 15.55864080031428
$ws.Range("G19").Value = 57.40352889213637
$ws.Range("H19").Value = 21.20038191327659
$ws.Range("J19").Value = 9.064686202810497
$ws.Range("K19").Value = 14.73398671389604
$ws.Range("M19").Value = 19.33133287336876
$ws.Range("N19").Value = 23.39424855584209

$ws.Range("B20").Value = 15.45707373729883
$ws.Range("C20").Value = 7.400104505630589
$ws.Range("D20").Value = 14.69826866717015
$ws.Range("E20").Value = 15.56416943092124
$ws.Range("G20").Value = 57.47808492647329
$ws.Range("H20").Value = 21.1977138480099
$ws.Range("J20").Value = 9.059419348583594
$ws.Range("K20").Value = 14.78212855058081
$ws.Range("M20").Value = 19.35148813036732
$ws.Range("N20").Value = 23.37767303910486

$ws.Range("B21").Value = 15.71106832571573
$ws.Range("C21").Value = 7.5572606229282
$ws.Range("D21").Value = 14.75253003999645
$ws.Range("E21").Value = 15.58594094766768
$ws.Range("G21").Value = 57.74251190409372
$ws.Range("H21").Value = 21.19260987642372
$ws.Range("J21").Value = 9.04268692534966
$ws.Range("K21").Value = 14.94678438067638
$ws.Range("M21").Value = 19.42348000519937
$ws.Range("N21").Value = 23.32422479470873

$ws.Range("B22").Value = 15.87816909717788
$ws.Range("C22").Value = 7.659901601658256
$ws.Range("D22").Value = 14.79041522446767
$ws.Range("E22").Value = 15.6025397631524
$ws.Range("G22").Value = 57.92556773771725
$ws.Range("H22").Value = 21.19215046183495
$ws.Range("J22").Value = 9.032467867024733
$ws.Range("K22").Value = 15.05640831616762
$ws.Range("M22").Value = 19.47367035616194
$ws.Range("N22").Value = 23.29096442667812

$ws.Range("B23").Value = 15.78890134656823
$ws.Range("C23").Value = 7.60513833418115
$ws.Range("D23").Value = 14.76997571601148
$ws.Range("E23").Value = 15.59346465446393
$ws.Range("G23").Value = 57.82694032846996
$ws.Range("H23").Value = 21.19213256280182
$ws.Range("J23").Value = 9.037856506108231
$ws.Range("K23").Value = 14.99772584097545
$ws.Range("M23").Value = 19.44659856794493
$ws.Range("N23").Value = 23.30856324616789

$ws.Range("B24").Value = 15.45294489863266
$ws.Range("C24").Value = 7.397537363593759
$ws.Range("D24").Value = 14.69742315298499
$ws.Range("E24").Value = 15.56385377219694
$ws.Range("G24").Value = 57.47393791530481
$ws.Range("H24").Value = 21.19784584163356
$ws.Range("J24").Value = 9.059704950424184
$ws.Range("K24").Value = 14.77947351052702
$ws.Range("M24").Value = 19.35036511420713
$ws.Range("N24").Value = 23.37857484610046

$ws.Range("B25").Value = 15.09709678104617
$ws.Range("C25").Value = 7.174531037634464
$ws.Range("D25").Value = 14.62976552137717
$ws.Range("E25").Value = 15.54216321542701
$ws.Range("G25").Value = 57.138017124954
$ws.Range("H25").Value = 21.21636058129284
$ws.Range("J25").Value = 9.086358548756357
$ws.Range("K25").Value = 14.55367932519197
$ws.Range("M25").Value = 19.26032447319684
$ws.Range("N25").Value = 23.46127385741722
